$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.094470262527466
$ws.Range("B1").Value = 5.52388334274292
$ws.Range("C1").Value = 2.472086906433105
$ws.Range("D1").Value = 1.643358111381531
$ws.Range("E1").Value = 1.648900032043457
